$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated roster table (header row A1:C1 is unchanged: "Oyuncu Adı", "Pozisyon", "Takım")
$data = @(
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Stephon Castle", "PG", "San Antonio Spurs"),
    @("Jalen Green", "PG", "Houston Rockets"),
    @("Brandin Podziemski", "PG,SG", "Golden State Warriors"),
    @("Dillon Brooks", "SG,SF,PF", "Houston Rockets"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Pascal Siakam", "SF,PF,C", "Indiana Pacers"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Russell Westbrook", "PG,SG", "Denver Nuggets"),
    @("Anthony Black", "PG,SG", "Orlando Magic"),
    @("Payton Pritchard", "PG,SG", "Boston Celtics"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Jalen Suggs", "PG,SG", "Orlando Magic"),
    @("Jakob Poeltl", "C", "Toronto Raptors")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
